$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.197.03'
$ws.Range("E2").Value = '  +1.03%  '
$ws.Range("D3").Value = '1.859.10'
$ws.Range("E3").Value = '  +1.25%  '
$ws.Range("E4").Value = '  +0.59%  '
$ws.Range("E5").Value = '  +3.39%  '
$ws.Range("D6").Value = "'0.623"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.60%  '
$ws.Range("E7").Value = '  +0.60%  '
$ws.Range("D8").Value = "'42.30"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.91%  '
$ws.Range("E9").Value = '  +0.78%  '
$ws.Range("E10").Value = '  +1.47%  '
$ws.Range("D11").Value = "'0.0989"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.30%  '
$ws.Range("D12").Value = '2.127.72'
$ws.Range("E12").Value = '  +1.34%  '
$ws.Range("D13").Value = "'11.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.51%  '
$ws.Range("D14").Value = '1.863.34'
$ws.Range("E14").Value = '  +1.56%  '
$ws.Range("E15").Value = '  +0.67%  '
$ws.Range("D16").Value = "'4.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.00%  '
$ws.Range("D17").Value = '35.159.53'
$ws.Range("E17").Value = '  +0.94%  '
$ws.Range("D18").Value = "'69.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.42%  '
$ws.Range("E19").Value = '  +1.20%  '
$ws.Range("D20").Value = "'240.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.14%  '
$ws.Range("D21").Value = "'12.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.66%  '
$ws.Range("E22").Value = '  +1.39%  '
$ws.Range("E23").Value = '  +0.52%  '
$ws.Range("E24").Value = '  +0.74%  '
$ws.Range("D25").Value = "'169.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.11%  '
$ws.Range("D26").Value = "'1.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +25.82%  '
$ws.Range("D27").Value = "'8.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.68%  '
$ws.Range("D28").Value = "'17.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.92%  '
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("E30").Value = '  +1.58%  '
$ws.Range("E31").Value = '  +0.56%  '
$ws.Range("E32").Value = '  +2.01%  '
$ws.Range("E33").Value = '  +27.32%  '
$ws.Range("E34").Value = '  +2.43%  '
$ws.Range("E35").Value = '  +10.68%  '
$ws.Range("D36").Value = "'0.818"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +17.26%  '
$ws.Range("D37").Value = "'1.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.52%  '
$ws.Range("E38").Value = '  +3.73%  '
$ws.Range("D39").Value = "'0.0202"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.29%  '
$ws.Range("D40").Value = "'90.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.58%  '
$ws.Range("D41").Value = '1.349.46'
$ws.Range("E41").Value = '  +0.73%  '
$ws.Range("D42").Value = "'0.0595"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +13.97%  '
$ws.Range("D43").Value = "'14.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.12%  '
$ws.Range("E44").Value = '  +2.56%  '
$ws.Range("D46").Value = "'12.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +43.88%  '
$ws.Range("E47").Value = '  -0.82%  '
$ws.Range("E48").Value = '  +4.88%  '
$ws.Range("D49").Value = '2.044.01'
$ws.Range("E49").Value = '  +1.52%  '
$ws.Range("D50").Value = "'0.0684"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.58%  '
$ws.Range("E51").Value = '  +0.57%  '
